$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table's first data row (2009年) is being dropped; every remaining year
# shifts up by one row, and a new 2021年 row is appended at the bottom.
$ws.Rows.Item(2).Delete()

# Append the new 2021年 data row at the bottom (now row 13).
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 25.8313
$ws.Range("C13").Value = 1.4215
$ws.Range("D13").Value = 22.931
$ws.Range("E13").Value = 55.7738
$ws.Range("F13").Value = 22.0862
$ws.Range("G13").Value = 134.1655
$ws.Range("H13").Value = 228.4567
$ws.Range("I13").Value = 8.4711
$ws.Range("J13").Value = 370.8813
$ws.Range("K13").Value = 3.1127
$ws.Range("L13").Value = 2277.8977
$ws.Range("M13").Value = 2235.343
$ws.Range("N13").Value = 453.9622
$ws.Range("O13").Value = 473.7532
$ws.Range("P13").Value = 467.1902
$ws.Range("Q13").Value = 10.0696
$ws.Range("R13").Value = 3.6113
$ws.Range("S13").Value = 11.6803
$ws.Range("T13").Value = 0.5344
$ws.Range("U13").Value = 221.296
$ws.Range("V13").Value = 30.7169

# Match the new label cell's style to the other year cells (bold, bordered, centered).
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A13").Value = "2021年"
